$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three CMS "Name" cells in column A to append their product URL,
# matching the fuller "Name - URL" convention already used for Drupal/DNN rows.
$ws.Range("A2").Value = "concrete5 - https://www.concrete5.org/"
$ws.Range("A5").Value = "Liferay Portal (ENTERPRISE SOLUTION ONLY) - https://www.liferay.com/product/features"
$ws.Range("A6").Value = "WordPress - https://en.wordpress.com/features/"

# Move the view: scroll so row 3 is the top visible row, and select A6.
$ws.Activate()
$ws.Range("A6").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
